$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B4").Value = 9.192599999999997
$ws.Range("B6").Value = 5.434500000000002
$ws.Range("B7").Value = 5.575000000000004
$ws.Range("E7").Value = 15.6235
$ws.Range("E12").Value = 17.6105
$ws.Range("E15").Value = 16.1207
$ws.Range("B16").Value = 7.314399999999999
$ws.Range("B20").Value = 9.511399999999989
$ws.Range("E20").Value = 15.99729999999999
$ws.Range("E21").Value = 16.81170000000001
$ws.Range("E22").Value = 16.57130000000001
$ws.Range("E23").Value = 16.13149999999999
$ws.Range("B28").Value = 5.549900000000004
$ws.Range("B29").Value = 5.208800000000004
$ws.Range("E29").Value = 17.25900000000001
$ws.Range("B32").Value = 7.507099999999995
$ws.Range("E34").Value = 16.99290000000001
$ws.Range("B40").Value = 9.198599999999994
$ws.Range("E42").Value = 16.61609999999999
$ws.Range("E43").Value = 17.2744
$ws.Range("E44").Value = 16.6145
$ws.Range("E45").Value = 16.4863
$ws.Range("B46").Value = 5.8993
$ws.Range("E46").Value = 16.7159
$ws.Range("E50").Value = 16.45329999999999
$ws.Range("B51").Value = 5.312399999999999
$ws.Range("E51").Value = 17.16670000000001
$ws.Range("B52").Value = 5.028900000000001
$ws.Range("B57").Value = 5.082899999999997
$ws.Range("B59").Value = 4.974799999999997
$ws.Range("B62").Value = 5.475000000000001
$ws.Range("B66").Value = 5.962599999999997
$ws.Range("E66").Value = 17.0501
$ws.Range("E67").Value = 17.07170000000002
$ws.Range("B73").Value = 8.4184
$ws.Range("B74").Value = 9.004299999999995
$ws.Range("E79").Value = 18.32580000000002
$ws.Range("E84").Value = 16.61139999999999
$ws.Range("B92").Value = 4.628499999999998
$ws.Range("E92").Value = 18.76140000000002
$ws.Range("E97").Value = 16.5606
$ws.Range("B100").Value = 5.300400000000001
